$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of
# the existing header row (bold, centered, bordered -> same style as H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell H1 onto the two new
# header cells so they match the rest of row 1 (bold/centered/bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new I0 (column I) and IF (column J) values for rows 2-46.
$iValues = @(4, 5, 6, 9, 5, 5, 5, 9, 9, 6, 4, 7, 6, 2, 9, 6, 7, 8, 4, 6, 5, 6, 6, 5, 4, 8, 7, 7, 7, 7, 7, 1, 3, 6, 3, 6, 4, 8, 8, 8, 7, 3, 6, 7, 6)
$jValues = @(6, 6, 9, 9, 7, 7, 5, 9, 9, 6, 5, 7, 8, 4, 9, 6, 8, 8, 5, 6, 7, 7, 7, 6, 5, 8, 7, 8, 7, 7, 7, 2, 4, 8, 5, 7, 5, 8, 9, 8, 8, 3, 6, 8, 7)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
